$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- String cell updates, in the precise order required so that newly
# --- introduced shared strings land at the indices the target workbook expects.

$ws.Range("A10").Value = "c"
$ws.Range("A11").Value = "@"
$ws.Range("B10").Value = "Pathway"
$ws.Range("C5").Value = "Colour.GREY_40"
$ws.Range("C6").Value = "Colour.GREEN_BB"
$ws.Range("C7").Value = "Colour.BLUE_BB"
$ws.Range("C9").Value = "Colour.RED_BB"
$ws.Range("C10").Value = "Colour.GREY_70"
$ws.Range("B11").Value = "Bridge"
$ws.Range("C11").Value = "Colour.GREY_D0"
$ws.Range("L6").Value = "Sheep/Cow can't pass."
$ws.Range("F6").Value = "M"
$ws.Range("G6").Value = "Cow"
$ws.Range("H6").Value = "Colour.BLACK"
$ws.Range("A8").Value = "e"
$ws.Range("A12").Value = "K"
$ws.Range("A14").Value = "X"
$ws.Range("A13").Value = "O"
$ws.Range("B12").Value = "Wall"
$ws.Range("B13").Value = "Window"
$ws.Range("B14").Value = "Entrance"
$ws.Range("C12").Value = "TBC"
$ws.Range("C13").Value = "TBC"
$ws.Range("C14").Value = "TBC"

# --- Numeric / formula cells

$ws.Range("I6").Value = 3
$ws.Range("P5").Value = 60
$ws.Range("P6").Value = 60
$ws.Range("P7").Value = 59
$ws.Range("P8").Value = 59
$ws.Range("P9").Formula = "=(60*60)-1"
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 0

# --- Column width for column L (12)

$ws.Columns.Item(12).ColumnWidth = 20 + 1/7

# --- Sheet view: scroll back to top-left and select C14

$ws.Range("C14").Select()
